$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.094.38"
$ws.Range("E2").Value = "  +5.21%  "
$ws.Range("D3").Value = "2.292.26"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'231.46"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'60.91"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +4.84%  "
$ws.Range("D10").Value = "'0.0945"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").Value = "'57.91"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "2.633.38"
$ws.Range("D14").Value = "'24.30"
$ws.Range("E14").Value = "  +8.39%  "
$ws.Range("D15").Value = "'15.77"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "'5.95"
$ws.Range("E16").Value = "  +5.75%  "
$ws.Range("D17").Value = "'0.817"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "2.295.45"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").Value = "44.016.14"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").Value = "'73.78"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "'6.26"
$ws.Range("E22").Value = "  +3.56%  "
$ws.Range("D23").Value = "'254.70"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +6.64%  "
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").Value = "'171.30"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("E29").Value = "  -3.50%  "
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'5.09"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'8.77"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "'0.000224"
$ws.Range("E43").Value = "  -13.50%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0969"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'98.93"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "'1.21"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "'10.39"
$ws.Range("E47").Value = "  +18.28%  "
$ws.Range("D48").Value = "'17.08"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "1.483.23"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'4.40"
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("E51").Value = "  +1.02%  "
